$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# are re-shuffled across rows 2-34 per the target revision.
$data = @(
    @{ Row = 2; D = 44179; J = 78; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 3; D = 44225; J = 56; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 4; D = 44223; J = 80; K = 2500; L = 3000; M = 2781; P = 927 }
    @{ Row = 5; D = 44557; J = 104; K = 2000; L = 2500; M = 2260; P = 753 }
    @{ Row = 6; D = 44389; J = 81; K = 2800; L = 3000; M = 2889; P = 963 }
    @{ Row = 7; D = 44291; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 8; D = 44292; J = 40; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 9; D = 44669; J = 92; K = 2500; L = 3000; M = 2755; P = 918 }
    @{ Row = 10; D = 44965; J = 87; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 11; D = 44537; J = 88; K = 2000; L = 2200; M = 2091; P = 697 }
    @{ Row = 12; D = 44756; J = 104; K = 2800; L = 3000; M = 2904; P = 968 }
    @{ Row = 13; D = 44967; J = 110; K = 3000; L = 3300; M = 3136; P = 1045 }
    @{ Row = 14; D = 44804; J = 85; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 15; D = 44187; J = 65; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 16; D = 44193; J = 70; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 17; D = 44992; J = 45; K = 4000; L = 4000; M = 4000; P = 1333 }
    @{ Row = 18; D = 44165; J = 68; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 19; D = 44224; J = 67; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 20; D = 44242; J = 95; K = 2500; L = 3000; M = 2737; P = 912 }
    @{ Row = 21; D = 44845; J = 80; K = 2500; L = 2500; M = 2500; P = 833 }
    @{ Row = 22; D = 44390; J = 50; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 23; D = 44260; J = 60; K = 3500; L = 3500; M = 3500; P = 1167 }
    @{ Row = 24; D = 44935; J = 78; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 25; D = 44937; J = 68; K = 3500; L = 3500; M = 3500; P = 1167 }
    @{ Row = 26; D = 44221; J = 50; K = 2500; L = 2500; M = 2500; P = 833 }
    @{ Row = 27; D = 44340; J = 54; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 28; D = 44536; J = 125; K = 2200; L = 2200; M = 2200; P = 733 }
    @{ Row = 29; D = 44627; J = 78; K = 3500; L = 3500; M = 3500; P = 1167 }
    @{ Row = 30; D = 44166; J = 45; K = 2500; L = 2500; M = 2500; P = 833 }
    @{ Row = 31; D = 44243; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 32; D = 44574; J = 50; K = 3000; L = 3000; M = 3000; P = 1000 }
    @{ Row = 33; D = 44559; J = 68; K = 2000; L = 2000; M = 2000; P = 667 }
    @{ Row = 34; D = 44222; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 16).Value = $item.P
}
